# Insert a new record row at row 194 in the "Albahaca" data table.
# This shifts the existing rows 194-237 down to 195-238 (row 237's data now
# lands on row 238), and the newly inserted row 194 is populated with a new
# weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 194, pushing 194..237 -> 195..238
$ws.Rows.Item(194).Insert()

# Populate the new row 194 with the new record's data
$ws.Range("A194").Value = 10
$ws.Range("B194").Value = "Vega Modelo de Temuco"
$ws.Range("C194").Value = "La Araucanía"
$ws.Range("D194").Value = 44722
$ws.Range("E194").Value = 9
$ws.Range("F194").Value = 100112052
$ws.Range("G194").Value = "Albahaca"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 30
$ws.Range("K194").Value = 5000
$ws.Range("L194").Value = 5000
$ws.Range("M194").Value = 5000
$ws.Range("N194").Value = "$/paquete"
$ws.Range("O194").Value = "Región de Arica y Parinacota"
$ws.Range("P194").Value = 5000
$ws.Range("Q194").Value = 1
$ws.Range("R194").Value = "Hortaliza"

# Make sure the date cell keeps the date number format (style index 2 in
# styles.xml), matching the rest of column D.
$ws.Range("D194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
